$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.444.62"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "'1.878.66"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'0.7199"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").Value = "'240.36"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.07821"
$ws.Range("E8").Value = "  -2.22%  "

$ws.Range("D9").Value = "'0.3117"
$ws.Range("E9").Value = "  +2.92%  "

$ws.Range("D10").Value = "'25.07"
$ws.Range("E10").Value = "  +6.89%  "

$ws.Range("D11").Value = "'0.08249"
$ws.Range("E11").Value = "  +0.58%  "

$ws.Range("D12").Value = "'1.891.36"
$ws.Range("E12").Value = "  +5.72%  "

$ws.Range("D13").Value = "'0.7287"
$ws.Range("E13").Value = "  +3.62%  "

$ws.Range("D14").Value = "'5.295"
$ws.Range("E14").Value = "  +2.25%  "

$ws.Range("D15").Value = "'91.45"
$ws.Range("E15").Value = "  +2.06%  "

$ws.Range("D16").Value = "'29.582.81"
$ws.Range("E16").Value = "  +2.32%  "

$ws.Range("D17").Value = "'5.952"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "'246.63"
$ws.Range("E18").Value = "  +3.77%  "

$ws.Range("D19").Value = "'0.000007893"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").Value = "'13.32"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").Value = "'0.9992"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'7.967"
$ws.Range("E22").Value = "  +6.71%  "

$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'0.1577"
$ws.Range("E24").Value = "  +9.40%  "

$ws.Range("D25").Value = "'163.89"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").Value = "'9.048"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").Value = "'18.34"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("E28").Value = "  -4.20%  "

$ws.Range("D29").Value = "'1.483"
$ws.Range("E29").Value = "  +0.36%  "

$ws.Range("D30").Value = "'4.389"
$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").Value = "'4.150"
$ws.Range("E31").Value = "  +3.27%  "

$ws.Range("D32").Value = "'0.05282"
$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("D33").Value = "'1.948"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").Value = "'1.203"
$ws.Range("E34").Value = "  +3.77%  "

$ws.Range("D35").Value = "'0.7225"
$ws.Range("E35").Value = "  +1.59%  "

$ws.Range("D36").Value = "'2.676"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "'0.01865"
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").Value = "'1.233.59"
$ws.Range("E38").Value = "  +8.59%  "

$ws.Range("D39").Value = "'2.724"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'0.9070"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").Value = "'73.87"
$ws.Range("E41").Value = "  +4.90%  "

$ws.Range("D42").Value = "'6.113"
$ws.Range("E42").Value = "  +3.93%  "

$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'103.48"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.766"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000120"
$ws.Range("E47").Value = "  +0.67%  "

$ws.Range("D48").Value = "'2.918"
$ws.Range("E48").Value = "  +12.71%  "

$ws.Range("D49").Value = "'0.4335"
$ws.Range("E49").Value = "  +1.93%  "

$ws.Range("D50").Value = "'9.289"
$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("D51").Value = "'7.088"
$ws.Range("E51").Value = "  +1.98%  "
